$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.041.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.831.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6322'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9992'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07510'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2938'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07700'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.831.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.994'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6711'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009597'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.080'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.054.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9986'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.166'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1410'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.499'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.067'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05392'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.856'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7444'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.87%  '
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.651'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.241.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.15%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01798'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.746'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.658'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9014'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9996'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.979.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("E45").Value = '  +3.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5107'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4062'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.004'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05774'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.758'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.28%  '
